$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Try"
$ws.Range("B18").Value = "Try"
$ws.Range("B18").Select()
